# Automatische test-sync: 2025-08-14 22:00:50
# Appends a new log row (row 35) to the "Logs" sheet, extends the
# conditional-formatting ranges that cover the data rows, and bumps the
# "Intern verzoek / Actie voor medewerker" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 35

$ws.Cells.Item($newRow, 1).Value = "Interne taak"
$ws.Cells.Item($newRow, 2).Value = "kwaliteit@testbedrijf123.nl"
$ws.Cells.Item($newRow, 3).Value = "Leg dit even neer bij Koen."
$ws.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@testbedrijf123.nl."
$ws.Cells.Item($newRow, 6).Value = "2025-08-14 22:00:01"
$ws.Cells.Item($newRow, 7).Value = "Nee"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional formatting sqref for every column that covered
# rows 2-34 so it now covers rows 2-35, matching the grown data range.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $extendedRange = $ws.Range($col + "2:" + $col + "35")
    $fc = $extendedRange.FormatConditions.Item(1)
    $fc.ModifyAppliesToRange($extendedRange)
}

# Update the Dashboard summary count for "Intern verzoek / Actie voor
# medewerker" (26 -> 27) to reflect the newly logged row.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 27
